$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that needs to move
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row
# (rows 2 through 439).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
